$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 10374
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 6937
$ws.Range("F7").Value = 652
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 1307
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = 900
$ws.Range("F15").Value = 106
$ws.Range("F19").Value = 335
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 1016
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 16
$ws.Range("F26").Value = 719
$ws.Range("F27").Value = 2960
$ws.Range("F29").Value = 2013
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 991
$ws.Range("F34").Value = 36
$ws.Range("F36").Value = 3669
$ws.Range("F37").Value = 4366
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 913
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F8").Value = 35
$ws.Range("F9").Value = 67
$ws.Range("F10").Value = 0
$ws.Range("F13").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 6
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 12
$ws.Range("F23").Value = 0
$ws.Range("F26").Value = 108
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 0

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 95
$ws.Range("F5").Value = 36
$ws.Range("F6").Value = 12
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 12667
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1307
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 106
$ws.Range("F18").Value = 1420
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 5
$ws.Range("F27").Value = 249
$ws.Range("F28").Value = 2013
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 36
$ws.Range("F37").Value = 4366
$ws.Range("F38").Value = 49
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 641
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 26
$ws.Range("F46").Value = 4281
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 61
